$d = $word.ActiveDocument

function Get-ParaByText($text) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.TrimEnd("`r`a") -eq $text) {
            return $p
        }
    }
    return $null
}

# 1) Replace the "ANSWER QUESTIONS HERE" placeholder paragraph with the
#    team's three-paragraph retrospective answer.
$target = Get-ParaByText "ANSWER QUESTIONS HERE"
if ($target -ne $null) {
    $target.Range.Text = "In this phase the team worked frantically. Perhaps not surprisingly, a large difficulty with managing a team is just thinking of tasks for people to do. Fortunately, Jesse was able to quickly finish assigned tickets, and the game was implemented according to plan. `r" + `
                          "After re-adjusting our design documentation though, the inefficiencies of our game became clear. The team made good use of GitLabs to outline tasks and manage different development timelines, and most big features were reviewed in a merge-request and commented as review - overall we used Git effectively. It would have been nice to re-adjust our design, but we do not have time. Perhaps in the coming weeks we will be able to make the game better. `r" + `
                          "I think we could have managed time better, which will be made actionable by starting on the third phase as soon as this one ends. This was a streeful, but I suppose enjoyable experience. "
}

# 2) Replace Gideon's "ADD CONTRIBUTIONS" placeholder with his actual
#    contributions paragraph.
$target = Get-ParaByText "ADD CONTRIBUTIONS"
if ($target -ne $null) {
    $target.Range.Text = "Made and assigned issues on GitLab for completion. Created Rooms and Map generator, NPC dialogue, Room descriptions and game layout. Reviewed merge requests and discussed new features in scrum-like meetings with Jesse after class. "
}

# 3) Move the "_GoBack" bookmark off the tail of Jess Huss's contributions
#    paragraph and collapse the three empty paragraphs that trail the
#    "Justin Creig" heading down to a single empty paragraph, which is where
#    the bookmark ends up living.
if ($d.Bookmarks.Exists("_GoBack")) {
    $null = $d.Bookmarks.Item("_GoBack").Delete()
}

# Merge the first of the three empty trailing paragraphs into the next one,
# twice, so only one empty paragraph remains after "Justin Creig".
for ($n = 0; $n -lt 2; $n++) {
    $justinPara = Get-ParaByText "Justin Creig"
    $firstEmpty = $justinPara.Next()
    $mergeRange = $firstEmpty.Range.Duplicate
    $null = $mergeRange.MoveEnd(1, 1)
    $null = $mergeRange.Delete()
}

# Re-add the bookmark, collapsed, at the start of the one remaining empty
# paragraph following "Justin Creig".
$justinPara = Get-ParaByText "Justin Creig"
$finalEmpty = $justinPara.Next()
$bmRange = $finalEmpty.Range.Duplicate
$null = $bmRange.Collapse(1)
$null = $d.Bookmarks.Add("_GoBack", $bmRange)
